# Finalizado Análisis del programa
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the input values in the "Análisis" block (rows 29-33). The
# dependent SUM formulas in A34 (=SUM(A20:A33)) and A62
# (=SUM(A18,A34,A41,A51,A59)) recalculate automatically from these inputs.
$ws.Range("A29").Value = 2
$ws.Range("A30").Value = 1
$ws.Range("A31").Value = 1
$ws.Range("A32").Value = 1
$ws.Range("A33").Value = 28

# Make sure everything is recalculated.
$excel.CalculateFull()

# Move the selection to the totals cell, matching the saved view state.
$ws.Activate()
$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
